# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" (Total) sheet,
#    populated with the per-fund holding breakdown for the new quarter
#    (same layout/style as the existing quarterly sheets).
# 2. Prepend a new summary row for "2022-Q1" to the "总计" sheet, shifting
#    the existing quarterly summary rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q1" sheet, positioned immediately before 总计
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($total)
$q1.Name = "2022-Q1"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 0; $col -lt $headers.Length; $col++) {
    $cell = $q1.Cells.Item(1, 2 + $col)
    $cell.Value = $headers[$col]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$fundRows = @(
    @("161903", "万家行业优选混合 (LOF)", "112.51", "91.94", "5.57", "6.2668", 9),
    @("008120", "万家自主创新混合A", "31.86", "93.78", "7.67", "2.4437", 6),
    @("005311", "万家经济新动能混合A", "19.47", "93.80", "8.37", "1.6296", 3),
    @("506000", "南方科创板 3 年定期开放混合", "24.62", "96.87", "5.94", "1.4624", 1),
    @("001643", "汇丰晋信智造先锋股票A", "29.09", "92.99", "3.64", "1.0589", 9),
    @("005312", "万家经济新动能混合C", "6.77", "93.80", "8.37", "0.5666", 3),
    @("001644", "汇丰晋信智造先锋股票C", "10.91", "92.99", "3.64", "0.3971", 9),
    @("008633", "万家科技创新混合A", "3.75", "93.27", "8.25", "0.3094", 5),
    @("008121", "万家自主创新混合C", "2.80", "93.78", "7.67", "0.2148", 6),
    @("008634", "万家科技创新混合C", "1.18", "93.27", "8.25", "0.0974", 5),
    @("000965", "汇丰晋信新动力混合", "0.99", "91.64", "3.45", "0.0342", 5),
    @("360012", "光大保德信中小盘混合", "1.10", "85.68", "2.55", "0.0280", 10),
    @("540004", "汇丰晋信2026周期混合", "1.14", "31.29", "1.61", "0.0184", 8)
)

for ($i = 0; $i -lt $fundRows.Length; $i++) {
    $r = 2 + $i
    $row = $fundRows[$i]

    $aCell = $q1.Cells.Item($r, 1)
    $aCell.Value = $i
    $aCell.Font.Bold = $true
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160
    $aCell.Borders.LineStyle = 1

    $bCell = $q1.Cells.Item($r, 2)
    $bCell.NumberFormat = "@"
    $bCell.Value = $row[0]

    $cCell = $q1.Cells.Item($r, 3)
    $cCell.Value = $row[1]

    $dCell = $q1.Cells.Item($r, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $row[2]

    $eCell = $q1.Cells.Item($r, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $row[3]

    $fCell = $q1.Cells.Item($r, 6)
    $fCell.NumberFormat = "@"
    $fCell.Value = $row[4]

    $gCell = $q1.Cells.Item($r, 7)
    $gCell.NumberFormat = "@"
    $gCell.Value = $row[5]

    $hCell = $q1.Cells.Item($r, 8)
    $hCell.Value = $row[6]
}

# ---------------------------------------------------------------------
# Step 2: insert a new top data row in "总计" for 2022-Q1, pushing the
# existing rows (2021-Q4, 2021-Q3, ...) down by one.
#
# NB: worksheet handles returned by Worksheets.Item()/Add() track sheet
# *position*, not identity - inserting the new "2022-Q1" sheet before
# 总计 shifts 总计 into a new slot, so re-resolve it by name before using
# it again (the stale $total handle now refers to the 2022-Q1 sheet).
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 1).Font.Bold = $true
$total.Cells.Item(2, 1).HorizontalAlignment = -4108
$total.Cells.Item(2, 1).VerticalAlignment = -4160
$total.Cells.Item(2, 1).Borders.LineStyle = 1

$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 13
$total.Cells.Item(2, 4).Value = 14.53

for ($r = 3; $r -le 7; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}

Write-Output "edit applied"
